$d = $word.ActiveDocument

# 1. Add the new "Abstract Title" paragraph style (wdStyleTypeParagraph = 1)
$title = $d.Styles.Add("Abstract Title", 1)
$title.BaseStyle = "Normal"
$title.NextParagraphStyle = "Abstract"
$title.QuickStyle = $true

$title.ParagraphFormat.KeepWithNext = $true
$title.ParagraphFormat.KeepTogether = $true
$title.ParagraphFormat.Alignment = 1
$title.ParagraphFormat.SpaceAfter = 0
$title.ParagraphFormat.SpaceBefore = 15

$title.Font.Size = 10
$title.Font.SizeBi = 10
$title.Font.Bold = $true
$title.Font.Color = "345A8A"

# 2. Tighten the space-before on the existing "Abstract" style (300 -> 100 twips = 5pt)
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5
